{"js": "const body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\nconsole.log(body.text);\n", "ps1": "$d = $word.ActiveDocument\n$r = $d.Content\n$x = $r.XML()\nWrite-Host \"----XML----\"\nWrite-Host $x\n"}
